$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: establish style pattern (s=6 plain / s=7 shaded) for the 2 brand-new rows
# by copying the existing row 11 formatting down into rows 12 and 13.
$ws.Range("A11:AB11").Copy($ws.Range("A12:AB12"))
$ws.Range("A11:AB11").Copy($ws.Range("A13:AB13"))

# --- Step 2: rewrite every data cell (rows 3-13) with the final values.
# A scratch cell (pre-formatted as Text) + PasteSpecial(values-only, -4163) is used for the
# column Z cells so the numeric-looking location strings land as literal text without
# disturbing the destination cell's existing style (a plain .Value= assignment of a
# numeric-looking string would silently be re-typed as a number by Excel).
$scratch = $ws.Range("AZ1")
$scratch.NumberFormat = "@"

# Row 3: PDH-E1 ALPHA
$ws.Range("A3").Value = 'PDH-E1 ALPHA'
$ws.Range("B3").Value = 'CAGTAACCTTGCCGAAGAG'
$ws.Range("C3").Value = 57.85284621491923
$ws.Range("D3").Value = 1.418635688603434
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0.007
$ws.Range("G3").Value = 1.1
$ws.Range("H3").Value = 'GCTCCTTACCAAGTCTGAC'
$ws.Range("I3").Value = 57.15171809660262
$ws.Range("J3").Value = 1.414071377081587
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0.007
$ws.Range("M3").Value = 3.7
$ws.Range("N3").Value = 2.832707065685021
$ws.Range("O3").Value = 'CCGTGACCATGTCCATGCCCTCAGCA'
$ws.Range("P3").Value = 64.46575342767647
$ws.Range("Q3").Value = 6.534246572323525
$ws.Range("R3").Value = -0.8420321941375732
$ws.Range("S3").Value = 0.00000151
$ws.Range("T3").Value = 1.4
$ws.Range("U3").Value = 'CAGTAACCTTGCCGAAGAGCTCGCTCATAACAGCACGAGCAGAGACACCTTTGCTGAGGGCATGGACATGGTCACGGTAGGTACTAACGACAGAGTCAGACTTGGTAAGGAGCT'
$ws.Range("V3").Value = 113
$ws.Range("W3").Value = 52.63157894736842
$ws.Range("X3").Value = -8.800433158874512
$ws.Range("Y3").Value = 'CP002684.1'
$scratch.Value = '48586'
$scratch.Copy()
$ws.Range("Z3").PasteSpecial(-4163)
$ws.Range("AA3").Value = 'LR699770.1'
$ws.Range("AB3").Value = 47493

# Row 4: PDH-E1 ALPHA
$ws.Range("A4").Value = 'PDH-E1 ALPHA'
$ws.Range("B4").Value = 'CGAAGAGCTCGCTCATAAC'
$ws.Range("C4").Value = 58.11854570998298
$ws.Range("D4").Value = 1.684335183667191
$ws.Range("E4").Value = -0.8030736446380615
$ws.Range("F4").Value = 0.007
$ws.Range("G4").Value = 1.1
$ws.Range("H4").Value = 'ACCAAGTCTGACTCTGTCG'
$ws.Range("I4").Value = 58.65308273860558
$ws.Range("J4").Value = 2.218872212289785
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0.007
$ws.Range("M4").Value = 3.7
$ws.Range("N4").Value = 3.903207395956976
$ws.Range("O4").Value = 'CCGTGACCATGTCCATGCCCTCAGCA'
$ws.Range("P4").Value = 64.46575342767647
$ws.Range("Q4").Value = 6.534246572323525
$ws.Range("R4").Value = -0.8420321941375732
$ws.Range("S4").Value = 0.00000151
$ws.Range("T4").Value = 1.4
$ws.Range("U4").Value = 'CGAAGAGCTCGCTCATAACAGCACGAGCAGAGACACCTTTGCTGAGGGCATGGACATGGTCACGGTAGGTACTAACGACAGAGTCAGACTTGGTA'
$ws.Range("V4").Value = 94
$ws.Range("W4").Value = 52.63157894736842
$ws.Range("X4").Value = -8.800433158874512
$ws.Range("Y4").Value = 'CP002684.1'
$scratch.Value = '48598'
$scratch.Copy()
$ws.Range("Z4").PasteSpecial(-4163)
$ws.Range("AA4").Value = 'LR699770.1'
$ws.Range("AB4").Value = 47505

# Row 5: CYP703A2
$ws.Range("A5").Value = 'CYP703A2'
$ws.Range("B5").Value = 'GAACGATCCCGATACCATC'
$ws.Range("C5").Value = 57.18013716608459
$ws.Range("D5").Value = 1.385652307599617
$ws.Range("E5").Value = -0.8097348213195801
$ws.Range("F5").Value = 0.007
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = 'CACATCCATACGCTAGGTG'
$ws.Range("I5").Value = 57.40147050241876
$ws.Range("J5").Value = 1.164318971265448
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0.007
$ws.Range("M5").Value = 3.7
$ws.Range("N5").Value = 2.549971278865065
$ws.Range("O5").Value = 'TCGAGACCCAAAACACTCGCCGCAGT'
$ws.Range("P5").Value = 64.30135325099201
$ws.Range("Q5").Value = 6.698646749007992
$ws.Range("R5").Value = -1.419864654541016
$ws.Range("S5").Value = 0.00000151
$ws.Range("T5").Value = 1.4
$ws.Range("U5").Value = 'GAACGATCCCGATACCATCCGTGAGATTCTTTTGCGGCAGGACGATGTTTTTTCATCGAGACCCAAAACACTCGCCGCAGTCCACCTAGCGTATGGATGTGG'
$ws.Range("V5").Value = 101
$ws.Range("W5").Value = 51.9607843137255
$ws.Range("X5").Value = -5.947404861450195
$ws.Range("Y5").Value = 'CP002684.1'
$scratch.Value = '112527'
$scratch.Copy()
$ws.Range("Z5").PasteSpecial(-4163)
$ws.Range("AA5").Value = 'LR699770.1'
$ws.Range("AB5").Value = 111452

# Row 6: CYP703A2
$ws.Range("A6").Value = 'CYP703A2'
$ws.Range("B6").Value = 'CAACGAACGATCCCGATAC'
$ws.Range("C6").Value = 58.1158368244387
$ws.Range("D6").Value = 1.681626298122911
$ws.Range("E6").Value = -0.8097348213195801
$ws.Range("F6").Value = 0.007
$ws.Range("G6").Value = 3.7
$ws.Range("H6").Value = 'ATCCATACGCTAGGTGGAC'
$ws.Range("I6").Value = 58.25082796726878
$ws.Range("J6").Value = 1.816617440952985
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0.007
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 3.498243739075896
$ws.Range("O6").Value = 'TCATCGAGACCCAAAACACTCGCCGCA'
$ws.Range("P6").Value = 64.20487603053124
$ws.Range("Q6").Value = 7.795123969468762
$ws.Range("R6").Value = -1.419864654541016
$ws.Range("S6").Value = 0.00000054
$ws.Range("T6").Value = 0.042
$ws.Range("U6").Value = 'CAACGAACGATCCCGATACCATCCGTGAGATTCTTTTGCGGCAGGACGATGTTTTTTCATCGAGACCCAAAACACTCGCCGCAGTCCACCTAGCGTATGGATG'
$ws.Range("V6").Value = 102
$ws.Range("W6").Value = 51.45631067961165
$ws.Range("X6").Value = -5.947404861450195
$ws.Range("Y6").Value = 'CP002684.1'
$scratch.Value = '112523'
$scratch.Copy()
$ws.Range("Z6").PasteSpecial(-4163)
$ws.Range("AA6").Value = 'LR699770.1'
$ws.Range("AB6").Value = 111448

# Row 7: At17.1
$ws.Range("A7").Value = 'At17.1'
$ws.Range("B7").Value = 'GTCTTCCGATCTCTCCATAC'
$ws.Range("C7").Value = 56.72627428247569
$ws.Range("D7").Value = 2.023725717524314
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0.002
$ws.Range("G7").Value = 0.3
$ws.Range("H7").Value = 'CTCAGCTTTGGGAATCCTC'
$ws.Range("I7").Value = 57.50429449531657
$ws.Range("J7").Value = 1.07008396900078
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0.007
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 3.093809686525095
$ws.Range("O7").Value = 'CCTCCGTCGTCCATCCTCCTGGTGGT'
$ws.Range("P7").Value = 64.72158538224431
$ws.Range("Q7").Value = 6.278414617755686
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0.00000151
$ws.Range("T7").Value = 0.4
$ws.Range("U7").Value = 'GTCTTCCGATCTCTCCATACTCAAACCACCAGGAGGATGGACGACGGAGGAAGAAGAAGGAGCGGTGAGGAGGGTGAAGGTGAGGATTCCCAAAGCTGAGC'
$ws.Range("V7").Value = 100
$ws.Range("W7").Value = 55.44554455445545
$ws.Range("X7").Value = -4.261404037475586
$ws.Range("Y7").Value = 'CP002685.1'
$scratch.Value = '164431'
$scratch.Copy()
$ws.Range("Z7").PasteSpecial(-4163)
$ws.Range("AA7").Value = 'LR699771.1'
$ws.Range("AB7").Value = 157366

# Row 8: At17.1
$ws.Range("A8").Value = 'At17.1'
$ws.Range("B8").Value = 'CATACTCAAACCACCAGGAG'
$ws.Range("C8").Value = 57.70701631138252
$ws.Range("D8").Value = 1.457016311382517
$ws.Range("E8").Value = -0.7797346115112305
$ws.Range("F8").Value = 0.002
$ws.Range("G8").Value = 1.1
$ws.Range("H8").Value = 'GGTGGCTCCTTCTTTAACG'
$ws.Range("I8").Value = 58.14276553363504
$ws.Range("J8").Value = 1.708555007319251
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0.007
$ws.Range("M8").Value = 3.7
$ws.Range("N8").Value = 3.165571318701769
$ws.Range("O8").Value = 'TCCTCACCTTCACCCTCCTCACCGCT'
$ws.Range("P8").Value = 64.17251465903956
$ws.Range("Q8").Value = 6.827485340960436
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0.00000151
$ws.Range("T8").Value = 0.033
$ws.Range("U8").Value = 'CATACTCAAACCACCAGGAGGATGGACGACGGAGGAAGAAGAAGGAGCGGTGAGGAGGGTGAAGGTGAGGATTCCCAAAGCTGAGCTTGAAAAGCTCGTTAAAGAAGGAGCCACCG'
$ws.Range("V8").Value = 115
$ws.Range("W8").Value = 53.44827586206896
$ws.Range("X8").Value = -5.429211139678955
$ws.Range("Y8").Value = 'CP002685.1'
$scratch.Value = '164446'
$scratch.Copy()
$ws.Range("Z8").PasteSpecial(-4163)
$ws.Range("AA8").Value = 'LR699771.1'
$ws.Range("AB8").Value = 157381

# Row 9: MTO1
$ws.Range("A9").Value = 'MTO1'
$ws.Range("B9").Value = 'CTGGTGGATCTAGGAGGTAC'
$ws.Range("C9").Value = 58.11536079407438
$ws.Range("D9").Value = 1.865360794074377
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0.002
$ws.Range("G9").Value = 0.3
$ws.Range("H9").Value = 'GGGATCAGGGAGAAGATAGG'
$ws.Range("I9").Value = 58.1693460818891
$ws.Range("J9").Value = 1.919346081889103
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0.002
$ws.Range("M9").Value = 1.1
$ws.Range("N9").Value = 3.784706875963479
$ws.Range("O9").Value = 'CCGATGGGAGCCTCACTGTTCATGCCG'
$ws.Range("P9").Value = 64.81452393084095
$ws.Range("Q9").Value = 7.185476069159051
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = 0.00000054
$ws.Range("T9").Value = 0.51
$ws.Range("U9").Value = 'CTGGTGGATCTAGGAGGTACCGGCATGAACAGTGAGGCTCCCATCGGAGCTCAAGAAAGAATGTTTGGAATCCGTCAGCTGTACACTACCTATCTTCTCCCTGATCCCC'
$ws.Range("V9").Value = 108
$ws.Range("W9").Value = 52.29357798165137
$ws.Range("X9").Value = -3.461651563644409
$ws.Range("Y9").Value = 'CP002686.1'
$scratch.Value = '41300'
$scratch.Copy()
$ws.Range("Z9").PasteSpecial(-4163)
$ws.Range("AA9").Value = 'LR699772.1'
$ws.Range("AB9").Value = 32887

# Row 10: MTO1
$ws.Range("A10").Value = 'MTO1'
$ws.Range("B10").Value = 'CTGAATCTGGTGGATCTAGG'
$ws.Range("C10").Value = 56.59552606387956
$ws.Range("D10").Value = 2.154473936120439
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0.002
$ws.Range("G10").Value = 1.1
$ws.Range("H10").Value = 'CAGGGAGAAGATAGGTAGTG'
$ws.Range("I10").Value = 55.9357109944566
$ws.Range("J10").Value = 2.814289005543401
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0.002
$ws.Range("M10").Value = 3.7
$ws.Range("N10").Value = 4.96876294166384
$ws.Range("O10").Value = 'CCGATGGGAGCCTCACTGTTCATGCCG'
$ws.Range("P10").Value = 64.81452393084095
$ws.Range("Q10").Value = 7.185476069159051
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0.00000054
$ws.Range("T10").Value = 0.51
$ws.Range("U10").Value = 'CTGAATCTGGTGGATCTAGGAGGTACCGGCATGAACAGTGAGGCTCCCATCGGAGCTCAAGAAAGAATGTTTGGAATCCGTCAGCTGTACACTACCTATCTTCTCCCTGA'
$ws.Range("V10").Value = 109
$ws.Range("W10").Value = 50
$ws.Range("X10").Value = -3.461651563644409
$ws.Range("Y10").Value = 'CP002686.1'
$scratch.Value = '41294'
$scratch.Copy()
$ws.Range("Z10").PasteSpecial(-4163)
$ws.Range("AA10").Value = 'LR699772.1'
$ws.Range("AB10").Value = 32881

# Row 11: XSP1
$ws.Range("A11").Value = 'XSP1'
$ws.Range("B11").Value = 'CTACGACATGGACGACATC'
$ws.Range("C11").Value = 57.44369071181501
$ws.Range("D11").Value = 1.1220987618692
$ws.Range("E11").Value = -0.3097348213195801
$ws.Range("F11").Value = 0.007
$ws.Range("G11").Value = 1.1
$ws.Range("H11").Value = 'GACAATGGAGGAACAGCTC'
$ws.Range("I11").Value = 57.83082272420569
$ws.Range("J11").Value = 1.396612197889898
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0.007
$ws.Range("M11").Value = 1.1
$ws.Range("N11").Value = 2.518710959759098
$ws.Range("O11").Value = 'AGCTCCATTGGTGGGTACACGCTCCG'
$ws.Range("P11").Value = 64.22598535419934
$ws.Range("Q11").Value = 6.774014645800662
$ws.Range("R11").Value = -0.7341697216033936
$ws.Range("S11").Value = 0.00000151
$ws.Range("T11").Value = 1.4
$ws.Range("U11").Value = 'CTACGACATGGACGACATCTCCTATGTTCAGTTCTTGTGCGGCGAAGGCTACAACGCAACCACTCTAGCTCCATTGGTGGGTACACGCTCCGTGAGCTGTTCCTCCATTGTCC'
$ws.Range("V11").Value = 112
$ws.Range("W11").Value = 53.98230088495575
$ws.Range("X11").Value = -6.521465301513672
$ws.Range("Y11").Value = 'CP002687.1'
$scratch.Value = '96840'
$scratch.Copy()
$ws.Range("Z11").PasteSpecial(-4163)
$ws.Range("AA11").Value = 'LR699773.1'
$ws.Range("AB11").Value = 149749

# Row 12: XSP1
$ws.Range("A12").Value = 'XSP1'
$ws.Range("B12").Value = 'CTGGCTTAGTCTACGACATG'
$ws.Range("C12").Value = 57.6139686963578
$ws.Range("D12").Value = 1.363968696357801
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0.002
$ws.Range("G12").Value = 3.7
$ws.Range("H12").Value = 'AGGGACAATGGAGGAACAG'
$ws.Range("I12").Value = 58.60725372965351
$ws.Range("J12").Value = 2.173043203337716
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0.007
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 3.537011899695518
$ws.Range("O12").Value = 'AGCTCCATTGGTGGGTACACGCTCCG'
$ws.Range("P12").Value = 64.22598535419934
$ws.Range("Q12").Value = 6.774014645800662
$ws.Range("R12").Value = -0.7341697216033936
$ws.Range("S12").Value = 0.00000151
$ws.Range("T12").Value = 1.4
$ws.Range("U12").Value = 'CTGGCTTAGTCTACGACATGGACGACATCTCCTATGTTCAGTTCTTGTGCGGCGAAGGCTACAACGCAACCACTCTAGCTCCATTGGTGGGTACACGCTCCGTGAGCTGTTCCTCCATTGTCCCTG'
$ws.Range("V12").Value = 125
$ws.Range("W12").Value = 53.96825396825397
$ws.Range("X12").Value = -6.521465301513672
$ws.Range("Y12").Value = 'CP002687.1'
$scratch.Value = '96830'
$scratch.Copy()
$ws.Range("Z12").PasteSpecial(-4163)
$ws.Range("AA12").Value = 'LR699773.1'
$ws.Range("AB12").Value = 149739

# Row 13: XSP1
$ws.Range("A13").Value = 'XSP1'
$ws.Range("B13").Value = 'CCCAACAATCCAACTCACG'
$ws.Range("C13").Value = 58.74250943646473
$ws.Range("D13").Value = 2.308298910148938
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0.007
$ws.Range("G13").Value = 3.7
$ws.Range("H13").Value = 'CCACCACTTTGAAGCTTCTC'
$ws.Range("I13").Value = 58.87532172899461
$ws.Range("J13").Value = 2.625321728994606
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0.002
$ws.Range("M13").Value = 0.3
$ws.Range("N13").Value = 4.933620639143545
$ws.Range("O13").Value = 'CCGCCACCGTCCGAGCACCG'
$ws.Range("P13").Value = 64.7335341111355
$ws.Range("Q13").Value = 0.2664658888645022
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0.002
$ws.Range("T13").Value = 1.1
$ws.Range("U13").Value = 'CCCAACAATCCAACTCACGTTGAGATCCGCCAAAACGTCCACATTGGCTGTGTTCAGGCGGAGAGTCACCAACGTGGGACCACCGTCGTCGGTCTACACCGCCACCGTCCGAGCACCGAAAGGAGTAGAAATCACGGTGGAGCCACAGAGTTTGTCATTTTCAAAGGCTTCACAAAAGAGAAGCTTCAAAGTGGTGGT'
$ws.Range("V13").Value = 197
$ws.Range("W13").Value = 53.03030303030303
$ws.Range("X13").Value = -17.26604652404785
$ws.Range("Y13").Value = 'CP002687.1'
$scratch.Value = '96981'
$scratch.Copy()
$ws.Range("Z13").PasteSpecial(-4163)
$ws.Range("AA13").Value = 'LR699773.1'
$ws.Range("AB13").Value = 149890

# --- Step 3: clean up the scratch cell so it does not show up in the used range / dimension
$scratch.Clear()

Write-Output "edit complete"
